$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 'Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(3, 7).Value = 'Dr. Manar Montaser, Dr. Gehan Adel, Administrator, Dr. Alshimaa Atef'
$ws.Cells.Item(4, 7).Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Asmaa Reda, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Majorelle Magdy'
$ws.Cells.Item(6, 7).Value = 'Dr. Safa Hany, Dr. Sara Nabil'
$ws.Cells.Item(9, 7).Value = 'Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Marina Youhanna, Dr. Yasmeena Fattoh'
$ws.Cells.Item(12, 7).Value = 'Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan, Dr. Mona Ibrahim Hussein'
$ws.Cells.Item(17, 7).Value = 'Dr. Walaa Ghanima, Dr. Enas Omran'
$ws.Cells.Item(18, 7).Value = 'Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Cells.Item(19, 7).Value = 'Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Cells.Item(20, 7).Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Yasmin, Dr. Marina Sorial, Dr. Nardine, Dr. Remon'
$ws.Cells.Item(21, 7).Value = 'Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(22, 7).Value = 'Dr. Manar Montaser, Dr. Gehan Adel, Administrator, Dr. Alshimaa Atef'
$ws.Cells.Item(23, 7).Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Nourhan Mahmoud, Dr. Hanan Ragab, Dr. Asmaa Reda, Dr. Shimaa Ahmad Mekki, Dr. Heba Mahmoud Ali, Dr. Majorelle Magdy'
$ws.Cells.Item(24, 7).Value = 'Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Cells.Item(25, 7).Value = 'Dr. Nourhan Mohammad, Dr. Yasmin Tarek'
$ws.Cells.Item(28, 7).Value = 'Dr. Eman M. Abo-Sakaya, Dr. Basma Hamed, Dr. Nourhan Osama, Dr. Madeha Saeed, Dr. Dina Adel, Dr. Yasmeena Fattoh, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Sarah Abdelmohsen, Dr. Marwa Mustafa'
$ws.Cells.Item(29, 7).Value = 'Dr. Amira Ibrahim, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh'
$ws.Cells.Item(31, 7).Value = 'Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan, Dr. Mona Ibrahim Hussein'
$ws.Cells.Item(36, 7).Value = 'Dr. Walaa Ghanima, Dr. Enas Omran'
$ws.Cells.Item(37, 7).Value = 'Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Cells.Item(38, 7).Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Yasmin, Dr. Marina Sorial, Dr. Nardine, Dr. Remon'
$ws.Cells.Item(39, 7).Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Yasmin, Dr. Marina Sorial, Dr. Nardine, Dr. Remon'
$ws.Cells.Item(40, 7).Value = 'Dr. Heba Mahmoud Ali, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Veronia Rafat, Dr. Amira Sobhy, Dr. Rana Abo-Zaid, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(41, 7).Value = 'Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Cells.Item(42, 7).Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(43, 7).Value = 'Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Kerelos Zareef'
$ws.Cells.Item(44, 7).Value = 'Dr. Safa Hany, Dr. Sara Nabil'
$ws.Cells.Item(47, 7).Value = 'Dr. Maryam Ahmad, Dr. Nourhan Osama, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Amira Ibrahim'
$ws.Cells.Item(48, 7).Value = 'Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Merna Said, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen'
$ws.Cells.Item(49, 7).Value = 'Dr. Mariam Toma Gerges, Dr. Mohammad Safwat'
$ws.Cells.Item(50, 7).Value = 'Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan, Dr. Mona Ibrahim Hussein'
$ws.Cells.Item(54, 7).Value = 'Dr. Amr Saeed, Dr. Afaf Abdallah'
$ws.Cells.Item(56, 7).Value = 'Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Cells.Item(57, 7).Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Yasmin, Dr. Marina Sorial, Dr. Nardine, Dr. Remon'
$ws.Cells.Item(58, 7).Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Yasmin, Dr. Marina Sorial, Dr. Nardine, Dr. Remon'
$ws.Cells.Item(59, 7).Value = 'Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Asmaa Reda, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(60, 7).Value = 'Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Cells.Item(61, 7).Value = 'Dr. Asmaa Reda, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Majorelle Magdy'
$ws.Cells.Item(63, 7).Value = 'Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed'
$ws.Cells.Item(66, 7).Value = 'Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Marina Youhanna'
$ws.Cells.Item(67, 7).Value = 'Dr. Amira Ibrahim, Dr. Esraa Mostafa, Dr. Yasmeena Fattoh'
$ws.Cells.Item(71, 7).Value = 'Dr. Nouran Mahmoud, Dr. Sarah Mahdy'
$ws.Cells.Item(72, 7).Value = 'Dr. Nouran Mahmoud, Dr. Sarah Mahdy'
$ws.Cells.Item(75, 7).Value = 'Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Cells.Item(76, 7).Value = 'Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Cells.Item(77, 7).Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Yasmin, Dr. Marina Sorial, Dr. Nardine, Dr. Remon'
$ws.Cells.Item(78, 7).Value = 'Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Asmaa Reda, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(79, 7).Value = 'Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Cells.Item(80, 7).Value = 'Dr. Asmaa Reda, Dr. Shimaa Ahmad Mekki, Dr. Amira Sobhy, Dr. Nahla Nagiub, Dr. Majorelle Magdy'
$ws.Cells.Item(81, 7).Value = 'Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Abeer Ragab'
$ws.Cells.Item(82, 7).Value = 'Dr. Nourhan Mohammad, Dr. Yasmin Tarek'
$ws.Cells.Item(83, 7).Value = 'Dr. Aya Saeed, Dr. Amal Awwad, Dr. Safa Hany'
$ws.Cells.Item(85, 7).Value = 'Dr. Eman M. Abo-Sakaya, Dr. Madeha Saeed, Dr. Yasmeena Fattoh, Dr. Dina Adel, Dr. Amira Ibrahim, Dr. Marina Youhanna'
$ws.Cells.Item(86, 7).Value = 'Dr. Eman M. Abo-Sakaya, Dr. Amany Raafat, Dr. Maryam Ahmad, Dr. Merna Said, Dr. Yasmeena Fattoh, Dr. Sarah Abdelmohsen'
$ws.Cells.Item(88, 7).Value = 'Dr. Heba Al-Sayed Mohammad, Dr. Dalia Tarek Elwan, Dr. Mona Ibrahim Hussein'
$ws.Cells.Item(90, 7).Value = 'Dr. Nouran Mahmoud, Dr. Sarah Mahdy'
$ws.Cells.Item(91, 7).Value = 'Dr. Nouran Mahmoud, Dr. Sarah Mahdy'
$ws.Cells.Item(94, 7).Value = 'Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Cells.Item(95, 7).Value = 'Dr. Eman Samir Gabry, Dr. Neveen Nashaat, Dr. Marina Sorial, Dr. Yasmin, Dr. Wafaa Ebida'
$ws.Cells.Item(96, 7).Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Yasmin, Dr. Marina Sorial, Dr. Nardine, Dr. Remon'
$ws.Cells.Item(97, 7).Value = 'Dr. Nourhan Mahmoud, Dr. Nesma, Dr. Asmaa Reda, Dr. Mohammad El-Tanany, Dr. Amira Sobhy, Dr. Heba Mahmoud Ali, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(98, 7).Value = 'Dr. Hanan Ragab, Dr. Mohammad El-Tanany, Dr. Shimaa Ahmad Mekki, Dr. Hend Mahmoud, Dr. Amira Sobhy, Dr. Alshimaa Atef'
$ws.Cells.Item(99, 7).Value = 'Dr. Menna tuâ€™Allah Medhat, Dr. Eman Tantawi, Dr. Shimaa Ahmad Mekki, Dr. Alshimaa Atef, Dr. Servinaz Sayed Mohammad'
$ws.Cells.Item(100, 7).Value = 'Dr. Nada Mohammad, Dr. Amera Ahmad Saad, Dr. Lamiaa Ossama, Dr. Menna tu''Alllah Mohammad, Dr. Fatma Elhady, Dr. Abeer Ragab, Dr. Kerelos Zareef'
$ws.Cells.Item(101, 7).Value = 'Dr. Safa Hany, Dr. Amal Awwad, Dr. Aya Saeed'
$ws.Cells.Item(104, 7).Value = 'Dr. Maryam Ahmad, Dr. Nourhan Osama, Dr. Merna Said, Dr. Esraa Mostafa, Dr. Arwa Al-Sayed, Dr. Amira Ibrahim'
$ws.Cells.Item(113, 7).Value = 'Dr. Abdullah El-Agrody, Dr. Ola Abd Al-Fattah, Dr. Eman Samir Gabry, Dr. Wafaa Ebida'
$ws.Cells.Item(115, 7).Value = 'Dr. Monica, Dr. Neveen Nashaat, Dr. Marina Atef, Dr. Yasmin, Dr. Marina Sorial, Dr. Nardine, Dr. Remon'
